# Auto-generated scheduled market-data refresh for Tiamat_Profits workbook
# Updates currentAveragePrice* and Leve profit columns (H:N) for affected leve rows
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 257
$ws.Range("I4").Value = 106.333336
$ws.Range("J4").Value = 386.14285
$ws.Range("K4").Value = 106.333336
$ws.Range("L4").Value = 386.14285
$ws.Range("M4").Value = 7.666663999999997
$ws.Range("N4").Value = -614.14285

$ws.Range("H70").Value = 1633.8695
$ws.Range("I70").Value = 1425
$ws.Range("J70").Value = 1745.2667
$ws.Range("K70").Value = 4275
$ws.Range("L70").Value = 5235.800099999999
$ws.Range("M70").Value = -4005
$ws.Range("N70").Value = -5775.800099999999

$ws.Range("H73").Value = 1633.8695
$ws.Range("I73").Value = 1425
$ws.Range("J73").Value = 1745.2667
$ws.Range("K73").Value = 4275
$ws.Range("L73").Value = 5235.800099999999
$ws.Range("M73").Value = -3339
$ws.Range("N73").Value = -7107.800099999999

$ws.Range("H111").Value = 1308.8572
$ws.Range("I111").Value = 1110.3334
$ws.Range("J111").Value = 1666.2
$ws.Range("K111").Value = 3331.0002
$ws.Range("L111").Value = 4998.6
$ws.Range("M111").Value = -264.0001999999999
$ws.Range("N111").Value = -11132.6

$ws.Range("H116").Value = 5296.6665
$ws.Range("I116").Value = 6415.885
$ws.Range("J116").Value = 4357.968
$ws.Range("K116").Value = 6415.885
$ws.Range("L116").Value = 4357.968
$ws.Range("M116").Value = -2973.885
$ws.Range("N116").Value = -11241.968

$ws.Range("H129").Value = 976.7273
$ws.Range("J129").Value = 994.44446
$ws.Range("L129").Value = 2983.33338
$ws.Range("N129").Value = -12983.33338

$ws.Range("H132").Value = 192806.48
$ws.Range("I132").Value = 4158.844
$ws.Range("J132").Value = 1253949.5
$ws.Range("K132").Value = 12476.532
$ws.Range("L132").Value = 3761848.5
$ws.Range("M132").Value = -9946.531999999999
$ws.Range("N132").Value = -3766908.5

$ws.Range("H137").Value = 5743.1904
$ws.Range("I137").Value = 950.1667
$ws.Range("J137").Value = 12133.889
$ws.Range("K137").Value = 2850.5001
$ws.Range("L137").Value = 36401.667
$ws.Range("M137").Value = -300.5001000000002
$ws.Range("N137").Value = -41501.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 27982.07
$ws.Range("I32").Value = 18842.07
$ws.Range("J32").Value = 41534.484
$ws.Range("K32").Value = 18842.07
$ws.Range("L32").Value = 41534.484
$ws.Range("M32").Value = -18555.07
$ws.Range("N32").Value = -42108.484

$ws.Range("H122").Value = 852.375
$ws.Range("I122").Value = 755.6316
$ws.Range("J122").Value = 1220
$ws.Range("K122").Value = 2266.8948
$ws.Range("L122").Value = 3660
$ws.Range("M122").Value = 183.1052
$ws.Range("N122").Value = -8560

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2893.2856
$ws.Range("I22").Value = 2893.2856
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2893.2856
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -2720.2856
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35494.684
$ws.Range("I31").Value = 45972
$ws.Range("J31").Value = 19429.467
$ws.Range("K31").Value = 45972
$ws.Range("L31").Value = 19429.467
$ws.Range("M31").Value = -45677
$ws.Range("N31").Value = -20019.467

$ws.Range("H34").Value = 35494.684
$ws.Range("I34").Value = 45972
$ws.Range("J34").Value = 19429.467
$ws.Range("K34").Value = 45972
$ws.Range("L34").Value = 19429.467
$ws.Range("M34").Value = -45770
$ws.Range("N34").Value = -19833.467

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3156.8096
$ws.Range("I3").Value = 2385.8333
$ws.Range("J3").Value = 4184.778
$ws.Range("K3").Value = 7157.499899999999
$ws.Range("L3").Value = 12554.334
$ws.Range("M3").Value = -7045.499899999999
$ws.Range("N3").Value = -12778.334

$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").ClearContents()

$ws.Range("H105").Value = 11301.7
$ws.Range("J105").Value = 11301.7
$ws.Range("L105").Value = 33905.10000000001
$ws.Range("N105").Value = -39147.10000000001

$ws.Range("H113").Value = 589.6
$ws.Range("I113").Value = 616.3333
$ws.Range("J113").Value = 582.9167
$ws.Range("K113").Value = 1848.9999
$ws.Range("L113").Value = 1748.7501
$ws.Range("M113").Value = 321.0001
$ws.Range("N113").Value = -6088.7501

$ws.Range("H129").Value = 49367.668
$ws.Range("I129").Value = 696.3333
$ws.Range("K129").Value = 2088.9999
$ws.Range("M129").Value = 2911.0001

$ws.Range("H131").Value = 167525.17
$ws.Range("I131").Value = 478
$ws.Range("J131").Value = 182711.27
$ws.Range("K131").Value = 1434
$ws.Range("L131").Value = 548133.8099999999
$ws.Range("M131").Value = 3606
$ws.Range("N131").Value = -558213.8099999999

$ws.Range("H134").Value = 4263.121
$ws.Range("I134").Value = 4060.125
$ws.Range("J134").Value = 4454.1763
$ws.Range("K134").Value = 12180.375
$ws.Range("L134").Value = 13362.5289
$ws.Range("M134").Value = -7110.375
$ws.Range("N134").Value = -23502.5289

$ws.Range("H137").Value = 40951080
$ws.Range("I137").Value = 133334610
$ws.Range("J137").Value = 7956962.5
$ws.Range("K137").Value = 400003830
$ws.Range("L137").Value = 23870887.5
$ws.Range("M137").Value = -399998730
$ws.Range("N137").Value = -23881087.5

$ws.Range("H138").Value = 1941.0454
$ws.Range("I138").Value = 1202.7142
$ws.Range("J138").Value = 3233.125
$ws.Range("K138").Value = 3608.1426
$ws.Range("L138").Value = 9699.375
$ws.Range("M138").Value = 1531.8574
$ws.Range("N138").Value = -19979.375

$ws.Range("H139").Value = 3597.5386
$ws.Range("I139").Value = 2207.7222
$ws.Range("J139").Value = 4788.8096
$ws.Range("K139").Value = 6623.1666
$ws.Range("L139").Value = 14366.4288
$ws.Range("M139").Value = -1483.1666
$ws.Range("N139").Value = -24646.4288

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 23035.479
$ws.Range("I132").Value = 1955.1875
$ws.Range("J132").Value = 65196.062
$ws.Range("K132").Value = 5865.5625
$ws.Range("L132").Value = 195588.186
$ws.Range("M132").Value = -3335.5625
$ws.Range("N132").Value = -200648.186

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 729.24
$ws.Range("I22").Value = 642.2857
$ws.Range("J22").Value = 839.9091
$ws.Range("K22").Value = 642.2857
$ws.Range("L22").Value = 839.9091
$ws.Range("M22").Value = -347.2857
$ws.Range("N22").Value = -1429.9091

$ws.Range("H27").Value = 729.24
$ws.Range("I27").Value = 642.2857
$ws.Range("J27").Value = 839.9091
$ws.Range("K27").Value = 642.2857
$ws.Range("L27").Value = 839.9091
$ws.Range("M27").Value = -535.2857
$ws.Range("N27").Value = -1053.9091

$ws.Range("H46").Value = 2036
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2036
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 2036
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -2412

$ws.Range("H55").Value = 239.72223
$ws.Range("I55").Value = 117.76923
$ws.Range("J55").Value = 556.8
$ws.Range("K55").Value = 117.76923
$ws.Range("L55").Value = 556.8
$ws.Range("M55").Value = 55.23077000000001
$ws.Range("N55").Value = -902.8

$ws.Range("H139").Value = 49224.645
$ws.Range("J139").Value = 49224.645
$ws.Range("L139").Value = 49224.645
$ws.Range("N139").Value = -59504.645

